$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Text "2024-02-28 Wednesday" "2024-02-29 Thursday"

# Table cell contents (ordered so that a replacement's new text never
# collides with an old-text search performed later)
Replace-Text "88×13=1144" "80×53=4240"
Replace-Text "97×96=9312" "88×43=3784"
Replace-Text "96×25=2400" "36×95=3420"
Replace-Text "16×25=400" "37×61=2257"
Replace-Text "30×82=2460" "19×13=247"
Replace-Text "83×85=7055" "70×99=6930"
Replace-Text "27×75=2025" "59×83=4897"
Replace-Text "58×32=1856" "26×19=494"
Replace-Text "88×48=4224" "20×69=1380"
Replace-Text "91×18=1638" "26×89=2314"
Replace-Text "21×25=525" "25×58=1450"
Replace-Text "46×97=4462" "25×81=2025"
Replace-Text "28×89=2492" "91×18=1638"
Replace-Text "82×70=5740" "96×17=1632"
Replace-Text "92×94=8648" "39×22=858"
Replace-Text "52×65=3380" "32×78=2496"
Replace-Text "94×50=4700" "45×56=2520"
Replace-Text "42×94=3948" "24×69=1656"
Replace-Text "26×51=1326" "72×29=2088"
Replace-Text "53×92=4876" "15×81=1215"
Replace-Text "23×46=1058" "25×99=2475"
Replace-Text "56×18=1008" "73×58=4234"
Replace-Text "62×38=2356" "40×21=840"
Replace-Text "23×16=368" "56×28=1568"
Replace-Text "99×75=7425" "87×11=957"
